$d = $word.ActiveDocument

function Set-RangeXml($range, $bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 4) "HIVE TEAM: QUALITY ASSURANCE" block: drop the "emelia" / "auditor" /
#        "Nitego" / "Senior QA Tester" paragraphs that followed "Release
#        Coordinator" entirely (they sat between "Release Coordinator" and
#        "clot").
$pStart = $d.Content.Find.Execute("emelia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rStart = $d.Content
$rStart.Find.Execute("emelia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $rStart.Paragraphs(1).Range.Start

$rEnd = $d.Content
$rEnd.Find.Execute("Senior QA Tester", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $rEnd.Paragraphs(1).Range.End

$d.Range($startPos, $endPos).Delete()

# --- 2) "HIVE TEAM: OUTREACH" block: drop the "Outreach Support" /
#        "Carlos Santiago" / "Outreach Support" / "Emilio" paragraphs that
#        followed "Semptly" entirely (they sat between "Semptly" and the
#        "Outreach Support" paragraph that precedes "misachasu").
$rStart2 = $d.Content
$rStart2.Find.Execute("Semptly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$semptlyEnd = $rStart2.Paragraphs(1).Range.End

$rEnd2 = $d.Content
$rEnd2.Find.Execute("Emilio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$emilioEnd = $rEnd2.Paragraphs(1).Range.End

$d.Range($semptlyEnd, $emilioEnd).Delete()

# --- 1) Replace "Mark Hakkarinen" / "Outreach Ambassador" / 'Editor of
#        "Your Week in SmartCash".' / "Email" (4 paragraphs) with the new
#        "LilyDaVine" heading + a single "Outreach Support" paragraph.
$rName = $d.Content
$rName.Find.Execute("Mark Hakkarinen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nameRange = $rName.Paragraphs(1).Range

$nameXml = '<w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="105" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr><w:t>LilyDaVine</w:t></w:r></w:p></w:body>'
Set-RangeXml $nameRange $nameXml

$rEmail = $d.Content
$rEmail.Find.Execute("Email", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$emailEnd = $rEmail.Paragraphs(1).Range.End

$rAmbassador = $d.Content
$rAmbassador.Find.Execute("Outreach Ambassador", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ambassadorStart = $rAmbassador.Paragraphs(1).Range.Start

$supportRange = $d.Range($ambassadorStart, $emailEnd)
$supportXml = '<w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="3B3B3B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Outreach Support</w:t></w:r></w:p></w:body>'
Set-RangeXml $supportRange $supportXml
